$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1651785714285714
$ws.Range("C2").Value = 0.6205357142857143
$ws.Range("J2").Value = 0.01785714285714286
$ws.Range("P2").Value = 0.1383928571428572
$ws.Range("S2").Value = 0.05803571428571429
# Row 3
$ws.Range("B3").Value = 0.01388888888888889
$ws.Range("C3").Value = 0.04166666666666666
$ws.Range("J3").Value = 0.03472222222222222
$ws.Range("P3").Value = 0.7569444444444444
$ws.Range("S3").Value = 0.1527777777777778
# Row 4
$ws.Range("J4").Value = 0.1081081081081081
$ws.Range("P4").Value = 0.6756756756756757
$ws.Range("S4").Value = 0.2162162162162162
# Row 5
$ws.Range("P5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.6666666666666666
# Row 6
$ws.Range("B6").Value = 0.05586592178770949
$ws.Range("F6").Value = 0.05586592178770949
$ws.Range("J6").Value = 0.2290502793296089
$ws.Range("O6").Value = 0.0111731843575419
$ws.Range("Q6").Value = 0.1005586592178771
$ws.Range("R6").Value = 0.0782122905027933
$ws.Range("S6").Value = 0.4692737430167598
# Row 7
$ws.Range("B7").Value = 0.1359223300970874
$ws.Range("D7").Value = 0.02912621359223301
$ws.Range("F7").Value = 0.009708737864077669
$ws.Range("J7").Value = 0.1650485436893204
$ws.Range("O7").Value = 0.04854368932038835
$ws.Range("Q7").Value = 0.1067961165048544
$ws.Range("R7").Value = 0.04854368932038835
$ws.Range("S7").Value = 0.4563106796116505
# Row 8
$ws.Range("B8").Value = 0.06188118811881188
$ws.Range("D8").Value = 0.01237623762376238
$ws.Range("E8").Value = 0.002475247524752475
$ws.Range("F8").Value = 0.06683168316831684
$ws.Range("J8").Value = 0.146039603960396
$ws.Range("O8").Value = 0.03465346534653466
$ws.Range("Q8").Value = 0.1584158415841584
$ws.Range("R8").Value = 0.1138613861386139
$ws.Range("S8").Value = 0.4034653465346535
# Row 9
$ws.Range("B9").Value = 0.09375
$ws.Range("D9").Value = 0.00390625
$ws.Range("F9").Value = 0.05859375
$ws.Range("J9").Value = 0.15234375
$ws.Range("O9").Value = 0.0390625
$ws.Range("Q9").Value = 0.17578125
$ws.Range("R9").Value = 0.09765625
$ws.Range("S9").Value = 0.37890625
# Row 10
$ws.Range("B10").Value = 0.09356223175965665
$ws.Range("D10").Value = 0.02489270386266094
$ws.Range("E10").Value = 0.001716738197424893
$ws.Range("F10").Value = 0.06266094420600858
$ws.Range("J10").Value = 0.1047210300429185
$ws.Range("O10").Value = 0.02317596566523605
$ws.Range("Q10").Value = 0.192274678111588
$ws.Range("R10").Value = 0.1012875536480687
$ws.Range("S10").Value = 0.3957081545064378
# Row 11
$ws.Range("G11").Value = 0.09523809523809523
$ws.Range("J11").Value = 0.1292517006802721
$ws.Range("K11").Value = 0.1224489795918367
$ws.Range("L11").Value = 0.6530612244897959
# Row 12
$ws.Range("G12").Value = 0.7156862745098039
$ws.Range("J12").Value = 0.2156862745098039
$ws.Range("K12").Value = 0.0196078431372549
$ws.Range("L12").Value = 0.04901960784313725
# Row 13
$ws.Range("G13").Value = 0.5294117647058824
$ws.Range("J13").Value = 0.4411764705882353
$ws.Range("S13").Value = 0.02941176470588235
# Row 15
$ws.Range("F15").Value = 0.004347826086956522
$ws.Range("H15").Value = 0.1608695652173913
$ws.Range("I15").Value = 0.08695652173913043
$ws.Range("J15").Value = 0.3608695652173913
$ws.Range("K15").Value = 0.03478260869565217
$ws.Range("M15").Value = 0.01304347826086956
$ws.Range("O15").Value = 0.03478260869565217
$ws.Range("S15").Value = 0.3043478260869565
# Row 16
$ws.Range("F16").Value = 0.0124223602484472
$ws.Range("H16").Value = 0.1677018633540373
$ws.Range("I16").Value = 0.1490683229813665
$ws.Range("J16").Value = 0.391304347826087
$ws.Range("K16").Value = 0.09316770186335403
$ws.Range("O16").Value = 0.04968944099378882
$ws.Range("S16").Value = 0.1366459627329193
# Row 17
$ws.Range("F17").Value = 0.02228412256267409
$ws.Range("H17").Value = 0.1615598885793872
$ws.Range("I17").Value = 0.1253481894150418
$ws.Range("J17").Value = 0.4484679665738162
$ws.Range("K17").Value = 0.06406685236768803
$ws.Range("M17").Value = 0.01671309192200557
$ws.Range("O17").Value = 0.07242339832869081
$ws.Range("S17").Value = 0.08913649025069638
# Row 18
$ws.Range("F18").Value = 0.02912621359223301
$ws.Range("H18").Value = 0.1650485436893204
$ws.Range("I18").Value = 0.1262135922330097
$ws.Range("J18").Value = 0.4320388349514563
$ws.Range("K18").Value = 0.04854368932038835
$ws.Range("M18").Value = 0.004854368932038835
$ws.Range("N18").Value = 0.004854368932038835
$ws.Range("O18").Value = 0.07281553398058252
$ws.Range("S18").Value = 0.116504854368932
# Row 19
$ws.Range("F19").Value = 0.01233480176211454
$ws.Range("H19").Value = 0.2237885462555066
$ws.Range("I19").Value = 0.1268722466960353
$ws.Range("J19").Value = 0.386784140969163
$ws.Range("K19").Value = 0.05991189427312775
$ws.Range("M19").Value = 0.02114537444933921
$ws.Range("N19").Value = 0.000881057268722467
$ws.Range("O19").Value = 0.07841409691629955
$ws.Range("S19").Value = 0.08986784140969163
